# Fill out the installation codes: add a new row for Walid's device entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "Walid"
$ws.Range("B15").Value = "RTX 4070"
$ws.Range("C15").Value = "16GB"

$ws.Range("C16").Select()
